# Refresh the cryptos price table (coinranking.com scrape) with the latest
# values pulled by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into Price (column D) while forcing it to stay plain
# text -- the sheet stores prices like "304.60" / "1.992.20" as literal
# strings (not numbers), but Excel.Range.Value auto-converts plain decimal
# strings into real numbers (dropping trailing zeros etc). Temporarily
# switching the cell to Text format forces the text type, then resetting
# the style back to "Normal" keeps the cell style byte-identical to before.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "42.568.29"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.287.98"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "304.60"
$ws.Range("E5").Value = "  +1.63%  "
Set-TextValue $ws.Range("D6") "95.36"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("E8").Value = "  +0.03%  "
Set-TextValue $ws.Range("D9") "0.495"
$ws.Range("E9").Value = "  -3.31%  "
Set-TextValue $ws.Range("D10") "34.86"
$ws.Range("E10").Value = "  -2.95%  "
Set-TextValue $ws.Range("D11") "0.0782"
$ws.Range("E11").Value = "  -0.67%  "
Set-TextValue $ws.Range("D12") "18.19"
$ws.Range("E12").Value = "  +3.49%  "
Set-TextValue $ws.Range("D14") "6.67"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "2.643.57"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "2.284.01"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "42.473.98"
$ws.Range("E18").Value = "  -0.96%  "
Set-TextValue $ws.Range("D19") "12.76"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  -2.38%  "
Set-TextValue $ws.Range("D21") "5.98"
$ws.Range("E21").Value = "  -1.98%  "
Set-TextValue $ws.Range("D22") "66.81"
$ws.Range("E22").Value = "  -2.88%  "
Set-TextValue $ws.Range("D23") "235.09"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").Value = "  -0.43%  "
Set-TextValue $ws.Range("D26") "2.43"
$ws.Range("E26").Value = "  +0.31%  "
Set-TextValue $ws.Range("D27") "4.02"
$ws.Range("E27").Value = "  +0.11%  "
Set-TextValue $ws.Range("D28") "24.88"
$ws.Range("E28").Value = "  +0.35%  "
Set-TextValue $ws.Range("D29") "165.54"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").Value = "  +0.50%  "
Set-TextValue $ws.Range("D31") "8.97"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D34") "4.93"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D35") "4.64"
$ws.Range("E35").Value = "  -1.68%  "
Set-TextValue $ws.Range("D36") "17.38"
$ws.Range("E36").Value = "  -1.93%  "
Set-TextValue $ws.Range("D37") "2.38"
$ws.Range("E37").Value = "  -0.90%  "
Set-TextValue $ws.Range("D38") "0.0683"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D40") "0.109"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D41") "1.73"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("D43").Value = "1.991.61"
$ws.Range("E43").Value = "  -0.46%  "
Set-TextValue $ws.Range("D44") "0.0276"
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D45") "9.96"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "17.84"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("E47").Value = "  -10.14%  "
$ws.Range("E48").Value = "  -1.93%  "
Set-TextValue $ws.Range("D49") "2.88"
$ws.Range("E49").Value = "  +8.46%  "
Set-TextValue $ws.Range("D50") "53.30"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "2.510.65"
$ws.Range("E51").Value = "  -0.35%  "
